$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old data rows 4-11 entirely (shrinks used range back down)
$ws.Range("A4:E11").ClearContents()

# Row 2: update client/plate/model/entry, clear exit (salida)
$ws.Range("A2").Value = "prueba"
$ws.Range("B2").Value = "prueba0"
$ws.Range("C2").Value = "honda"
$ws.Range("D2").Value = "2024-09-10 09:49 AM"
$ws.Range("E2").ClearContents()

# Row 3: clear client, update plate/model/entry, clear exit (salida)
$ws.Range("A3").ClearContents()
$ws.Range("B3").Value = "7163gmk"
$ws.Range("C3").Value = "honda"
$ws.Range("D3").Value = "2024-09-10 09:56 AM"
$ws.Range("E3").ClearContents()
